# Insert two new paragraphs ("Hola" and "No se entiende") right before
# the final (bookmark-only) paragraph of the document.
#
# Before:
#   P1: Cualquier Titulo
#   P2: Instrucciones de la actividad
#   P3: <empty, contains _GoBack bookmark>
#
# After:
#   P1: Cualquier Titulo
#   P2: Instrucciones de la actividad
#   P3: Hola
#   P4: No se entiende <empty, contains _GoBack bookmark>

$d = $word.ActiveDocument

# The last paragraph currently holds the _GoBack bookmark and no text.
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)

# Insert a brand-new paragraph "Hola" right before it.
$lastPara.Range.InsertBefore("Hola`r")

# The bookmark paragraph shifted down by one; grab it again and give it
# the "No se entiende" text (inserted at its very start, ahead of the
# bookmark start/end that already live there).
$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkPara.Range.InsertBefore("No se entiende")
